$p = $ppt.ActivePresentation

# The deck currently carries two themes: the active "Integral" theme used by
# the slide master (and therefore every slide), and a dormant "Office Theme"
# that is only linked from the notes master. Switching the presentation's
# design from "Integral" back to the (already-present) "Office Theme" swaps
# which of the two color palettes is active on the slide master - i.e. the
# master's theme color scheme becomes the Office Theme colors (and the old
# Integral colors end up the ones no longer driving the visible slides).

function RGBValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = RGBValue($officeThemeColors[$i - 1])
}
